$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.257266640663147
$ws.Range("B1").Value = 1.382025122642517
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.077202558517456
$ws.Range("E1").Value = 0.9087323546409607
